# Update both sheets (Stiffness = sheet1, Strength = sheet2) with new
# spatial-statistics values, and remove the now-unused 5th data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 2: "Strength"  (applied first so sheet1 ends up as the active /
# tabSelected tab, matching the original file's state)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Strength")

$ws2.Cells.Item(2,1).Value = 2
$ws2.Cells.Item(2,2).Value = 4
$ws2.Cells.Item(2,3).Value = 17
$ws2.Cells.Item(2,4).Value = 4
$ws2.Cells.Item(2,5).Value = 29.98
$ws2.Cells.Item(2,6).Value = 73.037999999999997

$ws2.Cells.Item(3,1).Value = 4
$ws2.Cells.Item(3,2).Value = 9
$ws2.Cells.Item(3,3).Value = 25
$ws2.Cells.Item(3,4).Value = 3
$ws2.Cells.Item(3,5).Value = 25.937000000000001
$ws2.Cells.Item(3,6).Value = 70.295000000000002

$ws2.Cells.Item(4,1).Value = 6
$ws2.Cells.Item(4,2).Value = 23
$ws2.Cells.Item(4,3).Value = 44
$ws2.Cells.Item(4,4).Value = 4
$ws2.Cells.Item(4,5).Value = 23.638999999999999
$ws2.Cells.Item(4,6).Value = 87.078000000000003

$ws2.Rows.Item(5).Delete()

$ws2.Range("B39").Select()

# ---------------------------------------------------------------------
# Sheet 1: "Stiffness"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Stiffness")

$ws1.Cells.Item(2,1).Value = 2
$ws1.Cells.Item(2,2).Value = 4
$ws1.Cells.Item(2,3).Value = 20
$ws1.Cells.Item(2,4).Value = 4
$ws1.Cells.Item(2,5).Value = 0.026100000000000002
$ws1.Cells.Item(2,6).Value = 0.126

$ws1.Cells.Item(3,1).Value = 4
$ws1.Cells.Item(3,2).Value = 9
$ws1.Cells.Item(3,3).Value = 26
$ws1.Cells.Item(3,4).Value = 2
$ws1.Cells.Item(3,5).Value = 0.20899999999999999
$ws1.Cells.Item(3,6).Value = 0.34699999999999998

$ws1.Cells.Item(4,1).Value = 6
$ws1.Cells.Item(4,2).Value = 23
$ws1.Cells.Item(4,3).Value = 14
$ws1.Cells.Item(4,4).Value = 3
$ws1.Cells.Item(4,5).Value = 177.53800000000001
$ws1.Cells.Item(4,6).Value = 551.30499999999995

# Row 5 no longer holds data - delete it so the used range shrinks to F4.
$ws1.Rows.Item(5).Delete()

$ws1.Range("A5").Select()
